$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table contained duplicate rows (id_heating_technology values 29, 210, 211, 33
# each appeared twice). Remove the duplicate occurrences, keeping the first one.
# These are Excel rows (1-indexed, row 1 = header): 12 (29), 14 (210), 16 (211), 22 (33).
# Delete from bottom to top so row numbers of the remaining rows-to-delete stay valid.
$rowsToDelete = @(22, 16, 14, 12)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete() | Out-Null
}

# Update the selected cell to B9 (matches the post-edit selection in the file)
$ws.Range("B9").Select() | Out-Null
